# Daily APR data refresh: append the latest scraped/snapshot row to Sheet1.
# New row 10: id=9, timestamp="2025-10-11T20:33", redemption_rate=1.7042856802003863
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "2025-10-11T20:33"
$ws.Range("C10").Value = 1.7042856802003863
